# Updated loading-percent results for the 380 kV case (rows 2-25, cols C-G/I-L/N-O).
# Columns H and M remain 0 and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 14, 15)

$data = @{
    2 = @(9.822975188313588, 8.382557364165249, 13.54097890451894, 37.33424484278362, 3.684782927488232, 28.99449201422622, 10.38122023858667, 17.53470595995335, 10.40610512897701, 18.17491623888462, 28.5868639829312)
    3 = @(9.779702990359894, 8.354175454606992, 13.53789774736287, 37.40571310126639, 3.68698184167519, 29.06349678656478, 10.40019537623019, 17.16509199144393, 10.42151919379589, 18.20924359476, 28.67259474696727)
    4 = @(9.754594240186213, 8.337728429668791, 13.53790823747973, 37.45820806963913, 3.688404093187684, 29.11200529675742, 10.4128981472917, 16.93594555692027, 10.43202919095558, 18.23216165333249, 28.7313746763993)
    5 = @(9.744737042435887, 8.33127683732034, 13.53839243528994, 37.48176154494237, 3.689001862299768, 29.13331372561631, 10.41833953523385, 16.84213479153787, 10.43657546209402, 18.24196463205231, 28.75686867913954)
    6 = @(9.743123112718152, 8.330220820799717, 13.53850186453122, 37.48580297890222, 3.689102221735029, 29.13694494409394, 10.41925908433402, 16.82653515485825, 10.43734628526461, 18.24362043724672, 28.7611948890968)
    7 = @(9.754459775268549, 8.337640400521581, 13.53791282240377, 37.45851697523123, 3.688412081183877, 29.11228643509386, 10.41297045865573, 16.93468197602691, 10.43208943667498, 18.23229198127715, 28.73171226391534)
    8 = @(9.807756456248255, 8.372570661390052, 13.53952267929756, 37.35709689474807, 3.685526181537941, 29.01700894063876, 10.38754476218997, 17.40779191889295, 10.41120311376709, 18.1863706871257, 28.61514736458965)
    9 = @(9.923509538118156, 8.44864020605859, 13.55769936182692, 37.22673323183734, 3.680436468497683, 28.87900610832028, 10.34601499378944, 18.31296140968827, 10.37852403491443, 18.11089267965676, 28.43543911175928)
    10 = @(10.01491003651341, 8.508862641004185, 13.58010211690467, 37.17293884025679, 3.677040578284438, 28.80754246124745, 10.32055835191111, 18.95772206088191, 10.35953712904453, 18.06427903018037, 28.33339986785056)
    11 = @(10.05775805556041, 8.53713652583478, 13.59223184738417, 37.15761475254154, 3.675569510329285, 28.78155760127, 10.31007025589733, 19.2453838232509, 10.35198463757825, 18.04498357645402, 28.29353072991332)
    12 = @(10.0741560295952, 8.547963714497078, 13.59710123603893, 37.15312865963802, 3.675022999518517, 28.77265769590366, 10.3062553513282, 19.35340628640979, 10.3492802220711, 18.0379507211161, 28.27937791705068)
    13 = @(10.07061693177053, 8.545626621719963, 13.59604028951372, 37.15403623959605, 3.67514023199384, 28.77453261795155, 10.30706999529135, 19.33018361450271, 10.34985575533726, 18.03945320014583, 28.2823839141112)
    14 = @(10.05910373170789, 8.538024904474023, 13.59262693681271, 37.15721928258385, 3.675524337406103, 28.78080655260352, 10.30975326275946, 19.25428964394726, 10.35175902879623, 18.04439949323781, 28.29234741847331)
    15 = @(10.05207370360026, 8.533384147768444, 13.59057203741853, 37.15934050559068, 3.675760985509562, 28.78477198200259, 10.31141723916183, 19.20768126736604, 10.3529450826612, 18.04746489325187, 28.29857347471572)
    16 = @(10.01213449783061, 8.507032056994602, 13.57934819819574, 37.17412451170588, 3.677138195326813, 28.80937208115586, 10.32126572238924, 18.93880066282558, 10.36005249301856, 18.06557839877322, 28.33613745631481)
    17 = @(9.987951067927577, 8.491086784587898, 13.57295743552612, 37.1855381114568, 3.678001917592906, 28.82613587210236, 10.3275869498716, 18.77233443220189, 10.36469018258864, 18.07717901191711, 28.36086132500795)
    18 = @(9.974161515585726, 8.48199852409563, 13.56946427861545, 37.19296385459406, 3.678505651648201, 28.83639197318197, 10.33132558279163, 18.67606170628003, 10.36745978741465, 18.08403113101572, 28.37569809844274)
    19 = @(9.969513546337295, 8.4789358321877, 13.56831300202824, 37.19562589532735, 3.678677401578481, 28.83996990930386, 10.33260909157892, 18.64337819138977, 10.36841508175893, 18.08638203336192, 28.38082732226991)
    20 = @(9.990513076650585, 8.492775633621971, 13.57361885907849, 37.18423399713914, 3.677909254605431, 28.82428777580661, 10.32690340444959, 18.79011018620272, 10.36418592582523, 18.07592550774068, 28.3581656239575)
    21 = @(10.0624808423686, 8.540254490489708, 13.59362204934279, 37.15624860012134, 3.675411230412837, 28.77893822192872, 10.30896087186124, 19.2766069501698, 10.35119577322868, 18.04293921821443, 28.28939523210069)
    22 = @(10.11051567676556, 8.571984118427476, 13.60830338986795, 37.14563402462242, 3.673840100460109, 28.75477928668412, 10.29814764553058, 19.58922305678942, 10.34361240508411, 18.02297705761661, 28.24995785887953)
    23 = @(10.08479053703453, 8.554987384915661, 13.6003214624461, 37.15059663016945, 3.674673034666506, 28.76717147025646, 10.30383542372234, 19.42289278915947, 10.34757699888594, 18.03348538971461, 28.2705013992576)
    24 = @(9.989354437203582, 8.492011858815868, 13.57331926557492, 37.18482089627521, 3.677951125207767, 28.82512137405907, 10.32721210987679, 18.78207553285703, 10.36441357839625, 18.07649164759493, 28.35938241094975)
    25 = @(9.891041670414403, 8.427279065413122, 13.55118455242691, 37.25463997622996, 3.681752780566412, 28.91109410368405, 10.35636046726722, 18.07120100249183, 10.3864807756313, 18.12975589873738, 28.47880149179923)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = $vals[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) rows of loading_percent results."
